# Applies the cryptos-list refresh described by the commit diff:
# price/volume(1h) updates for most rows, plus a re-sort that shifts the
# Bittensor/Hedera/Maker/Dai block (rows 35-38) and the VeChain/ApeXProtocol
# pair (rows 47-48).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value. Columns D (Price) and E (Volume) hold numeric-looking
# strings (e.g. "1.00", "3.259.09") that must stay text, so NumberFormat is
# forced to "@" before assignment and the style is reset to Normal afterwards
# so no stray number-format style gets attached to the cell.
$updates = [ordered]@{
    "D2" = "65.114.22"
    "E2" = "  -6.97%  "
    "D3" = "3.263.70"
    "E3" = "  -8.35%  "
    "D4" = "1.00"
    "E4" = "  +0.03%  "
    "D5" = "551.43"
    "E5" = "  -6.86%  "
    "D6" = "177.86"
    "E6" = "  -9.79%  "
    "E7" = "  +0.00%  "
    "E8" = "  -5.05%  "
    "D9" = "3.262.04"
    "E9" = "  -8.04%  "
    "D10" = "0.184"
    "E10" = "  -12.10%  "
    "D11" = "0.580"
    "E11" = "  -7.83%  "
    "D12" = "46.84"
    "E12" = "  -11.86%  "
    "D13" = "0.0000260"
    "E13" = "  -10.60%  "
    "D14" = "8.45"
    "E14" = "  -9.31%  "
    "D15" = "3.781.50"
    "E15" = "  -8.51%  "
    "D16" = "602.10"
    "E16" = "  -7.18%  "
    "D17" = "17.76"
    "E17" = "  -4.14%  "
    "D18" = "65.087.79"
    "E18" = "  -7.01%  "
    "E19" = "  -4.47%  "
    "D20" = "3.251.34"
    "E20" = "  -8.80%  "
    "D21" = "11.30"
    "E21" = "  -10.80%  "
    "D22" = "0.893"
    "E22" = "  -7.82%  "
    "D23" = "17.40"
    "E23" = "  -4.72%  "
    "D24" = "101.39"
    "E24" = "  -2.81%  "
    "D25" = "4.95"
    "E25" = "  -9.53%  "
    "D26" = "3.94"
    "E26" = "  -11.21%  "
    "D27" = "5.97"
    "E27" = "  -1.24%  "
    "D28" = "2.65"
    "E28" = "  -10.37%  "
    "D29" = "9.27"
    "E29" = "  -10.14%  "
    "D30" = "8.57"
    "E30" = "  -11.34%  "
    "D31" = "30.14"
    "E31" = "  -9.70%  "
    "D32" = "3.81"
    "E32" = "  -12.40%  "
    "D33" = "6.18"
    "E33" = "  -9.91%  "
    "D34" = "10.94"
    "E34" = "  -7.52%  "
    "B35" = "Bittensor"
    "C35" = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
    "D35" = "536.40"
    "E35" = "  +4.30%  "
    "B36" = "Hedera"
    "C36" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
    "D36" = "0.103"
    "E36" = "  -7.53%  "
    "B37" = "Maker"
    "C37" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
    "D37" = "3.729.63"
    "E37" = "  -0.45%  "
    "B38" = "Dai"
    "C38" = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
    "D38" = "1.00"
    "E38" = "  -0.04%  "
    "D39" = "55.97"
    "E39" = "  -9.60%  "
    "D40" = "3.41"
    "E40" = "  -7.69%  "
    "D41" = "0.0₃0697"
    "E41" = "  -14.78%  "
    "D42" = "2.65"
    "E42" = "  -11.27%  "
    "E43" = "  -7.74%  "
    "D44" = "31.68"
    "E44" = "  -9.99%  "
    "D45" = "0.335"
    "E45" = "  -10.62%  "
    "D46" = "3.16"
    "E46" = "  +17.06%  "
    "B47" = "VeChain"
    "C47" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D47" = "0.0406"
    "E47" = "  -11.43%  "
    "B48" = "ApeXProtocol"
    "C48" = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
    "D48" = "3.16"
    "E48" = "  -7.61%  "
    "E49" = "  -6.88%  "
    "D50" = "2.56"
    "E50" = "  -11.80%  "
    "D51" = "0.998"
    "E51" = "  -0.12%  "
}

foreach ($cellRef in $updates.Keys) {
    $col = $cellRef -replace '[0-9]+$', ''
    $range = $ws.Range($cellRef)
    if ($col -eq "D" -or $col -eq "E") {
        $range.NumberFormat = "@"
        $range.Value = $updates[$cellRef]
        $range.Style = "Normal"
    } else {
        $range.Value = $updates[$cellRef]
    }
}
